# Applies the "Added NCI SVN location for test data dump" edit to the
# 9581_Distribution_Success_Specimens_Derivatives_SpecimenArray test case.
#
# Strategy: the textual changes in the diff are small, surgical rewrites
# of a handful of paragraphs (merging runs that had been split apart by
# w:proofErr spell/grammar markers, and replacing the old single-line
# "Import dump located at ..." sentence with a multi-paragraph block that
# points at the new NCI SVN dump locations). We locate each paragraph by
# its distinctive text and overwrite it in one shot via Range.InsertXML,
# which lets us supply the exact target OOXML for the paragraph (so no
# stray proofErr markers get reintroduced the way they would via
# Find/Replace or Range.Text). w:lastRenderedPageBreak placement is a
# pagination side effect that the host recomputes on save, so it is left
# alone here.

$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [string]$MatchText,
        [string]$InnerXml
    )
    # NOTE: always call this with positional args - named args (-MatchText
    # "...") are not bound correctly by this host's PowerShell subset.

    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq $MatchText) {
            $target = $p
            break
        }
    }
    if ($null -eq $target) {
        throw "Paragraph not found for: $MatchText"
    }

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.Range.InsertXML($pkg)
}

# 1) "Select Test case ID 9581 with short title Distribution_..." - merge
#    the two runs that w:proofErr spellStart/spellEnd had split apart.
Set-ParagraphXml `
    "Select Test case ID 9581 with short title Distribution_Success_Specimens_Derivatives_SpecimenArray`r" `
    ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>' +
     '<w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr>' +
     '<w:r><w:t>Select Test case ID 9581 with short title Distribution_Success_Specimens_Derivatives_SpecimenArray</w:t></w:r></w:p>')

# 2) "Prerequisite: Import dump located at /files/caTissue/dump and deploy
#    application" -> trim the paragraph back to "Prerequisite: " and add a
#    new multi-paragraph block describing the latest dump locations.
# NOTE: Range.InsertXML fuses the very last paragraph mark of the supplied
# content into the destination's trailing mark (same as pasting - the
# final </w:p> never survives as its own paragraph). Since the diff wants
# a real trailing empty <w:p/> kept intact, tack on one extra throwaway
# empty paragraph at the end so the wanted <w:p/> isn't the one that is
# last/absorbed.
Set-ParagraphXml `
    "Prerequisite: Import dump located at /files/caTissue/dump and deploy application`r" `
    (
        '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Prerequisite:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Import latest dump located at </w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Oracle: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/Oracle</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>MySQL: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/MySQL and deploy application.</w:t></w:r></w:p>' +
        '<w:p/>' +
        '<w:p/>'
    )

# 3) "Specimens_ " + "collected" runs (split by proofErr spellStart/spellEnd
#    that also wrapped "_CAKUT") merge into a single "Specimens_ collected"
#    run; "_CAKUT" stays its own run, just losing the now-empty proofErr
#    wrapper.
Set-ParagraphXml `
    "Select Query title Specimens_ collected_CAKUT _protocol to execute.`r" `
    ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
     '<w:r><w:t xml:space="preserve">Select Query title </w:t></w:r>' +
     '<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Specimens_ collected</w:t></w:r>' +
     '<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>_CAKUT</w:t></w:r>' +
     '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
     '<w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>_protocol</w:t></w:r>' +
     '<w:r><w:t xml:space="preserve"> to execute.</w:t></w:r></w:p>')

# 4) "Click on " + "Save" + " button on the csv file." (split by proofErr
#    gramStart/gramEnd) merge into a single run.
Set-ParagraphXml `
    "Click on Save button on the csv file.`r" `
    ('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
     '<w:r><w:t>Click on Save button on the csv file.</w:t></w:r></w:p>')

# 5) "Specimen label, specimen type, tissue side, tissue site, pathological
#    status and protocol participant identifier." loses its proofErr
#    gramStart/gramEnd wrapper (no text change).
Set-ParagraphXml `
    "Specimen label, specimen type, tissue side, tissue site, pathological status and protocol participant identifier.`r" `
    '<w:p><w:r><w:t>Specimen label, specimen type, tissue side, tissue site, pathological status and protocol participant identifier.</w:t></w:r></w:p>'
